# Update "想去人数" (F column) figures across sheets, reflecting refreshed
# output generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 2458
$ws1.Range("F6").Value  = 1601
$ws1.Range("F9").Value  = 605
$ws1.Range("F10").Value = 3486
$ws1.Range("F15").Value = 879
$ws1.Range("F17").Value = 1210
$ws1.Range("F18").Value = 1776
$ws1.Range("F20").Value = 437
$ws1.Range("F21").Value = 1534
$ws1.Range("F23").Value = 1822
$ws1.Range("F25").Value = 4204

# --- Sheet: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value  = 63
$ws2.Range("F23").Value = 100
$ws2.Range("F36").Value = 420
$ws2.Range("F47").Value = 25
$ws2.Range("F48").Value = 25

# --- Sheet: 本地生活 (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value  = 2535
$ws3.Range("F5").Value  = 2545
$ws3.Range("F11").Value = 2941
$ws3.Range("F12").Value = 443
$ws3.Range("F13").Value = 783
$ws3.Range("F14").Value = 175

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 2535
$ws4.Range("F5").Value  = 2458
$ws4.Range("F7").Value  = 2941
$ws4.Range("F9").Value  = 783
$ws4.Range("F14").Value = 605
$ws4.Range("F20").Value = 879
$ws4.Range("F21").Value = 63
$ws4.Range("F23").Value = 1210
$ws4.Range("F30").Value = 1776
$ws4.Range("F31").Value = 437
$ws4.Range("F33").Value = 1534
$ws4.Range("F35").Value = 100
$ws4.Range("F36").Value = 100
$ws4.Range("F40").Value = 1823
$ws4.Range("F44").Value = 4204
$ws4.Range("F45").Value = 420
$ws4.Range("F51").Value = 25
